$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.425.89'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.701.21'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '''219.30'
$ws.Range("E5").Value = '  +0.47%  '
$ws.Range("D6").Value = '''0.5512'
$ws.Range("E6").Value = '  +5.17%  '
$ws.Range("D7").Value = '''1.010'
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''0.2745'
$ws.Range("E8").Value = '  +1.73%  '
$ws.Range("E9").Value = '  +1.00%  '
$ws.Range("D10").Value = '''22.08'
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("D11").Value = '''0.07700'
$ws.Range("E11").Value = '  +2.55%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").Value = '''4.553'
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").Value = '1.688.60'
$ws.Range("E13").Value = '  -0.25%  '
$ws.Range("D14").Value = '''0.5847'
$ws.Range("E14").Value = '  +0.96%  '
$ws.Range("D15").Value = '''0.000008395'
$ws.Range("E15").Value = '  -0.98%  '
$ws.Range("D16").Value = '''65.84'
$ws.Range("E16").Value = '  +2.38%  '
$ws.Range("D17").Value = '26.474.03'
$ws.Range("E17").Value = '  +0.39%  '
$ws.Range("D18").Value = '''4.956'
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("E19").Value = '  +0.15%  '
$ws.Range("E20").Value = '  +1.08%  '
$ws.Range("D21").Value = '''192.35'
$ws.Range("E21").Value = '  +1.95%  '
$ws.Range("E22").Value = '  +1.18%  '
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '''148.87'
$ws.Range("E24").Value = '  +2.94%  '
$ws.Range("D25").Value = '''0.1331'
$ws.Range("E25").Value = '  +8.18%  '
$ws.Range("D26").Value = '''7.931'
$ws.Range("E26").Value = '  +2.92%  '
$ws.Range("D27").Value = '''15.83'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '''0.06294'
$ws.Range("E28").Value = '  -5.17%  '
$ws.Range("D29").Value = '''1.386'
$ws.Range("E29").Value = '  +2.94%  '
$ws.Range("D31").Value = '''3.609'
$ws.Range("E31").Value = '  +1.06%  '
$ws.Range("D32").Value = '''3.617'
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("D33").Value = '''1.690'
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("D34").Value = '''1.046'
$ws.Range("E34").Value = '  +1.91%  '
$ws.Range("D35").Value = '''0.6187'
$ws.Range("E35").Value = '  -0.39%  '
$ws.Range("E36").Value = '  +0.46%  '
$ws.Range("D37").Value = '''2.765'
$ws.Range("E37").Value = '  +2.49%  '
$ws.Range("D38").Value = '''0.01649'
$ws.Range("E38").Value = '  +2.06%  '
$ws.Range("D39").Value = '1.121.54'
$ws.Range("E39").Value = '  +1.29%  '
$ws.Range("D40").Value = '''6.171'
$ws.Range("E40").Value = '  -3.09%  '
$ws.Range("D41").Value = '''0.8849'
$ws.Range("E41").Value = '  +0.44%  '
$ws.Range("E42").Value = '  +0.03%  '
$ws.Range("D43").Value = '''101.36'
$ws.Range("E43").Value = '  +0.36%  '
$ws.Range("D44").Value = '1.851.91'
$ws.Range("E44").Value = '  +1.06%  '
$ws.Range("D45").Value = '''57.87'
$ws.Range("E45").Value = '  +2.05%  '
$ws.Range("D46").Value = '''0.00000000108'
$ws.Range("E46").Value = '  -2.72%  '
$ws.Range("D47").Value = '''8.234'
$ws.Range("E47").Value = '  +1.02%  '
$ws.Range("D48").Value = '''1.006'
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").Value = '''0.05278'
$ws.Range("E49").Value = '  +0.16%  '
$ws.Range("D50").Value = '''6.147'
$ws.Range("E50").Value = '  +1.35%  '
$ws.Range("E51").Value = '  -0.01%  '
